$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'workout knee pad'
$ws.Cells.Item(2, 1).Value = 'girls knee pads volleyball'
$ws.Cells.Item(3, 1).Value = 'youth indoor volleyball'
$ws.Cells.Item(4, 1).Value = 'boys compressions'
$ws.Cells.Item(5, 1).Value = 'knee sleeve for bursitis'
$ws.Cells.Item(6, 1).Value = 'stretch marks men'
$ws.Cells.Item(7, 1).Value = 'mens below the knee shorts'
$ws.Cells.Item(8, 1).Value = 'sport tights men'
$ws.Cells.Item(9, 1).Value = 'joint compression'
$ws.Cells.Item(10, 1).Value = 'softball sweats'
$ws.Cells.Item(11, 1).Value = 'youth sport leggings for boys'
$ws.Cells.Item(12, 1).Value = 'padded compression knee sleeve'
$ws.Cells.Item(13, 1).Value = 'youth volleyball shorts'
$ws.Cells.Item(14, 1).Value = 'knee hockey'
$ws.Cells.Item(15, 1).Value = 'girls hockey compression shorts'
$ws.Cells.Item(16, 1).Value = 'football pads adult'
$ws.Cells.Item(17, 1).Value = 'padded football shorts'
$ws.Cells.Item(18, 1).Value = 'basketball shorts for men xxl'
$ws.Cells.Item(19, 1).Value = 'womens softball compression shorts'
$ws.Cells.Item(20, 1).Value = 'basketball leg sleeves youth boys'
$ws.Cells.Item(21, 1).Value = 'football youth tights'
$ws.Cells.Item(22, 1).Value = 'paintball protective gear pants'
$ws.Cells.Item(23, 1).Value = 'catcher pads'
$ws.Cells.Item(24, 1).Value = 'size chart for men'
$ws.Cells.Item(25, 1).Value = 'paintball shorts padded'
$ws.Cells.Item(26, 1).Value = 'men spandex leggings'
$ws.Cells.Item(27, 1).Value = 'calf compression sleeve for youth'
$ws.Cells.Item(28, 1).Value = 'mens 3/4 shorts'
$ws.Cells.Item(29, 1).Value = 'tendinitis knee'
$ws.Cells.Item(30, 1).Value = 'basketball hopes'
$ws.Cells.Item(31, 1).Value = 'mens volleyball shorts'
$ws.Cells.Item(32, 1).Value = 'cycling pants men padded'
$ws.Cells.Item(33, 1).Value = 'knee guards volleyball'
$ws.Cells.Item(34, 1).Value = 'silicone calf pads'
$ws.Cells.Item(35, 1).Value = 'work knees pads'
$ws.Cells.Item(36, 1).Value = 'youth padded knee sleeves'
$ws.Cells.Item(37, 1).Value = 'knee protector running'
$ws.Cells.Item(38, 1).Value = 'sliding shorts youth boys'
$ws.Cells.Item(39, 1).Value = 'padded leg sleeve'
$ws.Cells.Item(40, 1).Value = 'snowboarding pants for men'
$ws.Cells.Item(41, 1).Value = 'basketball sweat pants for men'
$ws.Cells.Item(42, 1).Value = 'black baseball pants youth boys'
$ws.Cells.Item(43, 1).Value = '1/2x28 thread protector'
$ws.Cells.Item(44, 1).Value = 'six one knee pads'
$ws.Cells.Item(45, 1).Value = 'outdoor knee pads'
$ws.Cells.Item(46, 1).Value = 'shorts with leggings men'
$ws.Cells.Item(47, 1).Value = 'youth baseball compression'
$ws.Cells.Item(48, 1).Value = 'hip guards for fall protection'
$ws.Cells.Item(49, 1).Value = 'mens compression yoga pants'
$ws.Cells.Item(50, 1).Value = 'basketball under shorts'
$ws.Cells.Item(51, 1).Value = 'capris spandex'
$ws.Cells.Item(52, 1).Value = 'men leg compression pants'
$ws.Cells.Item(53, 1).Value = 'youth compression padded shorts'
$ws.Cells.Item(54, 1).Value = 'youth volleyball knee pads for girls'
$ws.Cells.Item(55, 1).Value = 'baseball short pants'
$ws.Cells.Item(56, 1).Value = 'youth sports leggings boys'
$ws.Cells.Item(57, 1).Value = 'youth athletic tights boys'
$ws.Cells.Item(58, 1).Value = 'compression pant for men'
$ws.Cells.Item(59, 1).Value = 'working knee pad'
$ws.Cells.Item(60, 1).Value = 'mens long shorts below knee'
$ws.Cells.Item(61, 1).Value = 'black youth football pants'
$ws.Cells.Item(62, 1).Value = 'baseball pants youth xxl'
$ws.Cells.Item(63, 1).Value = 'girls compression leggings'
$ws.Cells.Item(64, 1).Value = 'easy knee pads'
$ws.Cells.Item(65, 1).Value = 'compression knee sleeve youth'
$ws.Cells.Item(66, 1).Value = 'basketball gear for boys youth'
$ws.Cells.Item(67, 1).Value = 'men leggings black'
$ws.Cells.Item(68, 1).Value = 'boy leggings for sports youth'
$ws.Cells.Item(69, 1).Value = 'leggings men compression'
$ws.Cells.Item(70, 1).Value = 'hockey pads youth'
$ws.Cells.Item(71, 1).Value = 'under knee pad sleeves'
$ws.Cells.Item(72, 1).Value = 'weightlifting floor'
$ws.Cells.Item(73, 1).Value = 'anti sweat pads'
$ws.Cells.Item(74, 1).Value = 'mens knee shorts'
$ws.Cells.Item(75, 1).Value = 'leggings for sports men'
$ws.Cells.Item(76, 1).Value = 'girls knee pads volleyball youth'
$ws.Cells.Item(77, 1).Value = 'compression knee sleeve with pad'
$ws.Cells.Item(78, 1).Value = 'running pad'
$ws.Cells.Item(79, 1).Value = 'baseball stretch bands'
$ws.Cells.Item(80, 1).Value = 'softball long pants'
$ws.Cells.Item(81, 1).Value = 'elastic waisted pants for men'
$ws.Cells.Item(82, 1).Value = 'mountain bike knee pads for men'
$ws.Cells.Item(83, 1).Value = 'yoga pant for men'
$ws.Cells.Item(84, 1).Value = 'compression basketball knee sleeve'
$ws.Cells.Item(85, 1).Value = 'compression knee sleeves with pads'
$ws.Cells.Item(86, 1).Value = 'knee work pad'
$ws.Cells.Item(87, 1).Value = 'running compression pants'
$ws.Cells.Item(88, 1).Value = 'professional knee pads construction'
$ws.Cells.Item(89, 1).Value = 'compression pants for girls'
$ws.Cells.Item(90, 1).Value = 'womans softball sliding shorts'
$ws.Cells.Item(91, 1).Value = 'compression calf sleeve youth'
$ws.Cells.Item(92, 1).Value = 'performance compression knee sleeve'
$ws.Cells.Item(93, 1).Value = 'youth wrestling kneepads'
$ws.Cells.Item(94, 1).Value = 'knee pads for work women'
$ws.Cells.Item(95, 1).Value = 'baseball pants short'
$ws.Cells.Item(96, 1).Value = 'men capri shorts'
$ws.Cells.Item(97, 1).Value = 'volleyball knee pads girls'
$ws.Cells.Item(98, 1).Value = 'sit pad hiking'
$ws.Cells.Item(99, 1).Value = 'tactical pants knee pads'
$ws.Cells.Item(100, 1).Value = 'training pants men'
